$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (2-4) -------------------------------------------------
$data = @(
    @{ A = "MCH229-1"; C = "(1981-1983) COMMUNITY ACTIVISIM"; D = "1981-1983"; E = "Series"; F = "1 Box"; G = "LOCATION: 24F | GRAP COUNT NUMER: NONE" },
    @{ A = "MCH229-2"; C = "(1958-1962) REGIONAL BRANCHES AND AFFILIATED ORGANISATIONS. PUBLICATIONS AND EDUCATIONAL MATERIALS."; D = "1958-1962"; E = "Series"; F = "1 Box"; G = "LOCATION: 24G | GRAP COUNT NUMER: NONE" },
    @{ A = "MCH229-3"; C = "(1946-1981) PUBLICATIONS AND EDUCATIONAL MATERIALS."; D = "1946-1981"; E = "Series"; F = "1 Box"; G = "LOCATION: 24G | GRAP COUNT NUMER: NONE" }
)

$row = 2
foreach ($rec in $data) {
    $ws.Range("A$row").Value = $rec.A
    $ws.Range("C$row").Value = $rec.C
    $ws.Range("D$row").Value = $rec.D
    $ws.Range("E$row").Value = $rec.E
    $ws.Range("F$row").Value = $rec.F
    $ws.Range("G$row").Value = $rec.G

    # Apply the data-row font to every populated cell plus the trailing
    # (empty) H cell so it still carries the style.
    $rng = $ws.Range("A$row" + ":" + "H$row")
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 10
    $rng.Font.ThemeColor = 1

    $row = $row + 1
}

# --- View: freeze the header row & select the new data range ------------
$ws.Range("A2:I4").Select()
$excel.ActiveWindow.FreezePanes = $true
